# Actualizacion Datos Personales 4 nov
# Updates the summary totals across the "Totales Plantel 1P", "Totales Plantel Final",
# "Reprobados por Grupo" and "Totales Grupos" sheets to reflect the latest grade data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Totales Plantel 1P" and "Totales Plantel Final" share identical values for
# rows 2-7 (columns D:M = grade distribution, Repro, Por_Repro) and receive
# the exact same updates.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Totales Plantel 1P", "Totales Plantel Final")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("J2").Value = 31
    $ws.Range("K2").Value = 112
    $ws.Range("L2").Value = 99
    $ws.Range("M2").Value = 46.92

    $ws.Range("G3").Value = 17
    $ws.Range("I3").Value = 15
    $ws.Range("J3").Value = 37

    $ws.Range("E4").Value = 11
    $ws.Range("F4").Value = 7
    $ws.Range("G4").Value = 8
    $ws.Range("I4").Value = 10
    $ws.Range("J4").Value = 28
    $ws.Range("K4").Value = 106
    $ws.Range("L4").Value = 82
    $ws.Range("M4").Value = 43.62

    $ws.Range("E5").Value = 14
    $ws.Range("F5").Value = 17
    $ws.Range("H5").Value = 12
    $ws.Range("I5").Value = 21

    $ws.Range("H6").Value = 14
    $ws.Range("J6").Value = 13

    $ws.Range("E7").Value = 21
    $ws.Range("F7").Value = 6
    $ws.Range("I7").Value = 15
    $ws.Range("J7").Value = 39
}

# ---------------------------------------------------------------------------
# "Reprobados por Grupo"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Reprobados por Grupo")

$ws4.Range("I3").Value = 6
$ws4.Range("J3").Value = 13
$ws4.Range("K3").Value = 18
$ws4.Range("L3").Value = 58.06

$ws4.Range("F11").Value = 3
$ws4.Range("H11").Value = 3
$ws4.Range("I11").Value = 8

$ws4.Range("D15").Value = 1
$ws4.Range("E15").Value = 2
$ws4.Range("G15").Value = 0
$ws4.Range("J15").Value = 24
$ws4.Range("K15").Value = 11
$ws4.Range("L15").Value = 31.43

$ws4.Range("G16").Value = 2
$ws4.Range("I16").Value = 6

$ws4.Range("D17").Value = 0
$ws4.Range("E17").Value = 1

$ws4.Range("F18").Value = 4
$ws4.Range("G18").Value = 6

$ws4.Range("G19").Value = 1
$ws4.Range("H19").Value = 2
$ws4.Range("J19").Value = 25
$ws4.Range("K19").Value = 9
$ws4.Range("L19").Value = 26.47

$ws4.Range("D20").Value = 4
$ws4.Range("E20").Value = 5
$ws4.Range("G20").Value = 3
$ws4.Range("H20").Value = 7

$ws4.Range("G25").Value = 6
$ws4.Range("H25").Value = 5

$ws4.Range("H29").Value = 3
$ws4.Range("I29").Value = 1

$ws4.Range("D33").Value = 5
$ws4.Range("E33").Value = 2
$ws4.Range("H33").Value = 2
$ws4.Range("I33").Value = 6

# ---------------------------------------------------------------------------
# "Totales Grupos"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Totales Grupos")

$ws5.Range("C3").Value = 13
$ws5.Range("D3").Value = 41.94
$ws5.Range("G3").Value = 13
$ws5.Range("H3").Value = 41.94

$ws5.Range("C15").Value = 24
$ws5.Range("D15").Value = 68.57
$ws5.Range("G15").Value = 24
$ws5.Range("H15").Value = 68.57

$ws5.Range("C19").Value = 25
$ws5.Range("D19").Value = 73.53
$ws5.Range("G19").Value = 25
$ws5.Range("H19").Value = 73.53
